# Updated symbol list on Thu Dec 22 13:44:47 UTC 2022 with GitHub Actions
#
# Refreshes the scraped coinranking.com price/volume snapshot: most rows only
# get a new Price (column D), a handful of rows (13-26) were re-ranked so
# their Coin/Link/Price/Volume columns shift to reflect the new ordering.
#
# Price cells are stored as text in the source sheet (the scraper writes raw
# strings), so each numeric-looking price is assigned with a leading
# apostrophe to force Excel to keep it as text instead of auto-converting it
# to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- simple price refreshes -------------------------------------------------
$ws.Cells.Item(2, 4).Value = "'" + "243.98"      # D2  BNB
$ws.Cells.Item(3, 4).Value = "'" + "22.28"       # D3  OKB
$ws.Cells.Item(4, 4).Value = "'" + "5.411"       # D4  HuobiToken
$ws.Cells.Item(4, 5).Value = "3HuobiTokenHT"     # E4  (Best-in-24h badge removed)
$ws.Cells.Item(5, 4).Value = "'" + "0.05758"     # D5  Cronos
$ws.Cells.Item(6, 4).Value = "'" + "3.431"       # D6  GateToken
$ws.Cells.Item(8, 4).Value = "'" + "0.8137"      # D8  MXToken
$ws.Cells.Item(9, 4).Value = "'" + "0.8735"      # D9  FTXToken
$ws.Cells.Item(10, 4).Value = "'" + "0.1442"     # D10 WazirX
$ws.Cells.Item(11, 4).Value = "'" + "0.07346"    # D11 MandalaExchangeToken
$ws.Cells.Item(12, 4).Value = "'" + "0.03035"    # D12 LiechtensteinCryptoassetsExchange

# --- rows 13-26 re-ranked: Coin / Link / Price / Volume(1h) all shift ------
$ws.Cells.Item(13, 2).Value = "BitrueCoin"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Cells.Item(13, 4).Value = "'" + "0.03107"
$ws.Cells.Item(13, 5).Value = "12BitrueCoinBTR"

$ws.Cells.Item(14, 2).Value = "BitMartToken"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Cells.Item(14, 4).Value = "'" + "0.09402"
$ws.Cells.Item(14, 5).Value = "13BitMartTokenBMX"

$ws.Cells.Item(15, 2).Value = "BitForexToken"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Cells.Item(15, 4).Value = "'" + "0.001592"
$ws.Cells.Item(15, 5).Value = "14BitForexTokenBF"

$ws.Cells.Item(16, 2).Value = "CoinExToken"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Cells.Item(16, 4).Value = "'" + "0.04830"
$ws.Cells.Item(16, 5).Value = "15CoinExTokenCET"

$ws.Cells.Item(17, 2).Value = "TigerCash"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Cells.Item(17, 4).Value = "'" + "0.006388"
$ws.Cells.Item(17, 5).Value = "16TigerCashTCH"

$ws.Cells.Item(18, 2).Value = "HotbitToken"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Cells.Item(18, 4).Value = "'" + "0.004135"
$ws.Cells.Item(18, 5).Value = "17HotbitTokenHTB"

$ws.Cells.Item(19, 2).Value = "BitKan"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Cells.Item(19, 4).Value = "'" + "0.0009976"
$ws.Cells.Item(19, 5).Value = "18BitKanKAN"

$ws.Cells.Item(20, 2).Value = "NitroEx"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Cells.Item(20, 4).Value = "'" + "0.0001501"
$ws.Cells.Item(20, 5).Value = "19NitroExNTX"

$ws.Cells.Item(21, 2).Value = "LEO"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(21, 4).Value = "'" + "3.724"
$ws.Cells.Item(21, 5).Value = "20LEOLEO"

$ws.Cells.Item(22, 2).Value = "BTSEToken"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Cells.Item(22, 4).Value = "'" + "2.197"
$ws.Cells.Item(22, 5).Value = "21BTSETokenBTSE"

$ws.Cells.Item(23, 2).Value = "One"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Cells.Item(23, 4).Value = "'" + "0.01098"
$ws.Cells.Item(23, 5).Value = "22OneONEBestin24h"

$ws.Cells.Item(24, 2).Value = "BitpandaEcosystemToken"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Cells.Item(24, 4).Value = "'" + "0.3275"
$ws.Cells.Item(24, 5).Value = "23BitpandaEcosystemTokenBEST"

$ws.Cells.Item(25, 2).Value = "ProBitToken"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Cells.Item(25, 4).Value = "'" + "0.1311"
$ws.Cells.Item(25, 5).Value = "24ProBitTokenPROB"

$ws.Cells.Item(26, 2).Value = "MCDex"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Cells.Item(26, 4).Value = "'" + "4.189"
$ws.Cells.Item(26, 5).Value = "25MCDexMCB"

# --- remaining simple price refreshes --------------------------------------
$ws.Cells.Item(27, 4).Value = "'" + "0.0003204"  # D27 UpBots
$ws.Cells.Item(40, 4).Value = "'" + "0.03882"    # D40 IDEX
$ws.Cells.Item(41, 4).Value = "'" + "0.006724"   # D41 KickToken
$ws.Cells.Item(42, 4).Value = "'" + "0.1069"     # D42 BKEXToken
$ws.Cells.Item(43, 4).Value = "'" + "0.002802"   # D43 CEJI
$ws.Cells.Item(44, 4).Value = "'" + "0.007455"   # D44 LocalTraders
$ws.Cells.Item(45, 4).Value = "'" + "0.00005593" # D45 CoinLion
$ws.Cells.Item(48, 4).Value = "'" + "0.1472"     # D48 BOLO
